$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "*maa://24633 (56.17), *maa://30515 (69.9), *maa://34787 (73.68), maa://39402 (92.19), ***maa://20792 (11.93), ***maa://29083 (27.78)"
$ws.Range("T2").Value = "maa://22742 (91.12), *maa://20791 (62.16)"
$ws.Range("T3").Value = "maa://24617 (89.83), **maa://20790 (43.48), ***maa://37170 (16.42), maa://45854 (84.0)"
$ws.Range("X4").Value = "**maa://32495 (48.7), ***maa://31785 (22.22), maa://43217 (90.91), ***maa://36683 (28.26)"
$ws.Range("D6").Value = "maa://42407 (96.43)"
$ws.Range("T7").Value = "maa://21291 (85.11)"
$ws.Range("A8").Value = "更新日期：2025.03.08 13:15:15"
$ws.Range("T9").Value = "**maa://22866 (30.19), maa://26222 (98.08)"
$ws.Range("X9").Value = "maa://26223 (97.96)"
$ws.Range("AB11").Value = "maa://29912 (97.26), maa://22516 (88.37), *maa://20794 (52.24)"
$ws.Range("AF12").Value = "*maa://28932 (77.85), *maa://20106 (63.96), *maa://22769 (64.29)"
$ws.Range("D13").Value = "maa://24999 (92.17), maa://36673 (93.33), maa://25001 (85.71)"
$ws.Range("X13").Value = "maa://34957 (82.05), **maa://22768 (50.0)"
$ws.Range("AB14").Value = "maa://22764 (97.22)"
$ws.Range("D16").Value = "maa://21441 (96.41), maa://36679 (94.55), maa://37650 (97.3)"
$ws.Range("L18").Value = "maa://22466 (90.53), *maa://22732 (51.09)"
$ws.Range("T19").Value = "maa://24386 (99.17)"
$ws.Range("AB19").Value = "*maa://30709 (65.91), *maa://36668 (57.5)"
$ws.Range("D20").Value = "maa://21432 (90.34), maa://25198 (93.64), *maa://20795 (50.77), maa://36680 (91.18)"
$ws.Range("L20").Value = "maa://41331 (85.44)"
$ws.Range("P20").Value = "maa://37442 (95.24)"
$ws.Range("L21").Value = "maa://31731 (96.23)"
$ws.Range("L23").Value = "maa://39756 (95.81), maa://39875 (94.52)"
$ws.Range("AF24").Value = "maa://22523 (85.64), maa://36672 (80.7), maa://29910 (93.22), **maa://21440 (35.71), *maa://45831 (80.0)"
$ws.Range("D25").Value = "maa://29753 (95.19)"
$ws.Range("X25").Value = "*maa://29890 (79.17)"
$ws.Range("AB25").Value = "maa://31215 (88.03), maa://24516 (80.22), maa://26001 (87.5)"
$ws.Range("D26").Value = "maa://41802 (95.0)"
$ws.Range("X28").Value = "maa://39929 (90.68), maa://41749 (90.91), ***maa://39723 (13.89)"
$ws.Range("AF28").Value = "maa://36660 (92.37), *maa://36701 (66.67)"
$ws.Range("L29").Value = "maa://28432 (93.55), maa://28440 (80.36), maa://31400 (98.82), *maa://28650 (71.43)"
$ws.Range("AB30").Value = "maa://42979 (97.04), maa://45822 (100.0), *maa://45045 (80.0)"
$ws.Range("T32").Value = "maa://42859 (96.0), maa://41108 (88.0), maa://41238 (97.14), maa://45523 (100.0)"
$ws.Range("L35").Value = "maa://41296 (96.34)"
$ws.Range("L37").Value = "maa://45718 (97.99), *maa://47069 (73.33), maa://45789 (100.0)"
$ws.Range("T39").Value = "maa://45788 (82.18), maa://47079 (94.12), *maa://45790 (76.92)"
$ws.Range("H43").Value = "maa://22525 (92.41), maa://21284 (85.71)"
$ws.Range("T44").Value = "maa://39366 (89.19)"
$ws.Range("T45").Value = "**maa://39364 (36.36)"
$ws.Range("H46").Value = "maa://35931 (91.95), maa://43901 (93.1)"
$ws.Range("H53").Value = "maa://32534 (94.12), **maa://32434 (33.33)"
$ws.Range("H55").Value = "maa://32532 (91.96)"
$ws.Range("H59").Value = "maa://31270 (94.66), maa://27746 (82.3)"
$ws.Range("H60").Value = "*maa://40438 (70.31)"
